# Tidsregistrering i PTE projektet - Nada Omer
# Commit: "Vaegt felt handterer komma nu"
#
# The existing log for 16-3-2017 (serial 42810) only had 3 entries (rows
# 47-49). A new entry was inserted right after the first one, the three
# original entries' dates were bumped to 17-3-2017 (serial 42811), a new
# 4th entry was appended for that day, two quiet days (18-3 and 19-3) were
# logged with just name/date, and 20-3 + 21-3 got their own entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 47 loses its old "Implementer OC7" entry - wipe the activity
#    data but keep the date (42810) and the cell formatting in place.
# ---------------------------------------------------------------------
$ws.Range("B47:I47").ClearContents()

# ---------------------------------------------------------------------
# 2. Rows 50, 51, 54, 55, 56, 57 are brand-new "full" entries (date,
#    name, role, activity, start/end time, duration). Clone the number
#    formatting from row 46 for every column that carries one (A, B, C,
#    G, H and I - E and F are always unstyled), then fill in the values.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","G","H","I")) {
    $ws.Range($col + "46").Copy()
    foreach ($r in @(50,51,54,55,56,57)) {
        $ws.Range($col + $r).PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------
# 3. Rows 52 & 53 are "quiet day" entries - only Date + Name, Role stays
#    blank. Clone formatting for A, B and C only.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C")) {
    $ws.Range($col + "46").Copy()
    foreach ($r in @(52,53)) {
        $ws.Range($col + $r).PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Fill in the values column-by-column (rather than row-by-row) so
#    that brand new shared-string entries get appended to the shared
#    string table in the same order the author typed them in: every
#    new Aktivitet (column F) first, then every new duration (col I).
# ---------------------------------------------------------------------

# --- Column A: Dato ---
$ws.Range("A48").Value = 42811
$ws.Range("A49").Value = 42811
$ws.Range("A50").Value = 42811
$ws.Range("A51").Value = 42811
$ws.Range("A52").Value = 42812
$ws.Range("A53").Value = 42813
$ws.Range("A54").Value = 42814
$ws.Range("A55").Value = 42815
$ws.Range("A56").Value = 42815
$ws.Range("A57").Value = 42815

# --- Column B: Navn ---
$ws.Range("B50").Value = "NO"
$ws.Range("B51").Value = "NO"
$ws.Range("B52").Value = "NO"
$ws.Range("B53").Value = "NO"
$ws.Range("B54").Value = "NO"
$ws.Range("B55").Value = "NO"
$ws.Range("B56").Value = "NO"
$ws.Range("B57").Value = "NO"

# --- Column C: Rolle (deltager) ---
$ws.Range("C50").Value = "Nada H. A. Omer"
$ws.Range("C51").Value = "Nada H. A. Omer"
$ws.Range("C54").Value = "Nada H. A. Omer"
$ws.Range("C55").Value = "Nada H. A. Omer"
$ws.Range("C56").Value = "Nada H. A. Omer"
$ws.Range("C57").Value = "Nada H. A. Omer"

# --- Column E: Rolle (gyldig rolle for aktiviteten) ---
$ws.Range("E50").Value = "Implementer"
$ws.Range("E51").Value = "Implementer"
$ws.Range("E54").Value = "Implementer"
$ws.Range("E55").Value = "Test Designer"
$ws.Range("E56").Value = "Reviewer"
$ws.Range("E57").Value = "Implementer"

# --- Column F: Aktivitet ---
$ws.Range("F48").Value = "Implementer OC7"
$ws.Range("F49").Value = "Implementer OC5"
$ws.Range("F50").Value = "Implementer  Junit test til OC7 + setFnNewton + setFtNewton"
$ws.Range("F51").Value = "Implementer  Junit test til OC5"
$ws.Range("F54").Value = "GUI Design"
$ws.Range("F55").Value = "Bruger Test"
$ws.Range("F56").Value = "Review OC15"
$ws.Range("F57").Value = "GUI Design - SigmaN + rettelser til udskrivning af mellemregninger"

# --- Column G: Starttid ---
$ws.Range("G48").Value = 0.33680555555555558
$ws.Range("G49").Value = 0.3888888888888889
$ws.Range("G50").Value = 0.4236111111111111
$ws.Range("G51").Value = 0.52083333333333337
$ws.Range("G54").Value = 0.33680555555555558
$ws.Range("G55").Value = 0.34375
$ws.Range("G56").Value = 0.42708333333333331
$ws.Range("G57").Value = 0.45833333333333331

# --- Column H: Sluttid ---
$ws.Range("H48").Value = 0.38541666666666669
$ws.Range("H49").Value = 0.41666666666666669
$ws.Range("H50").Value = 0.52083333333333337
$ws.Range("H51").Value = 0.58333333333333337
$ws.Range("H54").Value = 0.66666666666666663
$ws.Range("H55").Value = 0.42708333333333331
$ws.Range("H56").Value = 0.44097222222222227
$ws.Range("H57").Value = 0.5625

# --- Column I: Samlet tid brugt ---
$ws.Range("I48").Value = "1 time : 10 min."
$ws.Range("I49").Value = "0 time : 40 min."
$ws.Range("I50").Value = "1 time : 50 min."
$ws.Range("I51").Value = "1 time : 30 min."
$ws.Range("I54").Value = "7 time : 55 min."
$ws.Range("I55").Value = "2 time : 00 min."
$ws.Range("I56").Value = "0 time : 20 min."
$ws.Range("I57").Value = "2 time : 30 min."

# ---------------------------------------------------------------------
# 5. The two data-validation lists (Rolle, column C, and the hidden
#    E-column role helper) need to grow to cover the extra rows.
#    Delete + re-add (in their original order) so the sqref grows from
#    C3:C49 / E3:E124 to C3:C57 / E3:E126.
# ---------------------------------------------------------------------
$ws.Range("C3:C49").Validation.Delete()
$ws.Range("E3:E124").Validation.Delete()
$ws.Range("C3:C57").Validation.Add(3, 1, 1, "=Deltagere")
$ws.Range("E3:E126").Validation.Add(3, 1, 1, "=GyldigeRoller")

# ---------------------------------------------------------------------
# 6. Leave the selection on the last-edited cell, matching where the
#    author's cursor ended up.
# ---------------------------------------------------------------------
$ws.Range("I57").Select()
